$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# --- Remove the old hyperlinks before the row shift moves their anchors,
# so we can re-attach them to the correct (shifted) cells afterwards. ---
$hyperlinkTargets = @{
    "C87" = "https://www.edmundoptics.com/p/2x-mitutoyo-bd-plan-apo-objective/45425/"
    "C88" = "https://www.edmundoptics.com/p/5x-mitutoyo-bd-plan-apo-objective/45426/"
    "C89" = "https://www.edmundoptics.com/p/75x-mitutoyo-bd-plan-apo-objective/45427/"
    "C90" = "https://www.edmundoptics.com/p/10x-mitutoyo-bd-plan-apo-objective/45428/"
    "C95" = "https://www.edmundoptics.com/p/mt-1-accessory-tube-lens/11488/"
    "C78" = "https://astronomy-imaging-camera.com/product/efw-mini"
    "C91" = "https://shop.mitutoyo.ch/web/mitutoyo/en_CH/mitutoyo/05.04.02/G%20Plan%20APO%2020X-t3%2C5/PR/378-847/index.xhtml;jsessionid=CFD835C47751CA929FEA384183C0A21D"
}
foreach ($addr in $hyperlinkTargets.Keys) {
    $ws.Range($addr).Hyperlinks.Delete()
}

# Insert a new row at position 49, pushing everything else (including the
# hyperlinked cells above) down by one row.
$ws.Rows.Item(49).Insert()

# Populate the new row 49 with the 1:1 relay lens part.
$ws.Range("A49").Value = "AC254-100-A-ML"
$ws.Range("B49").Value = "Thorlabs"
$ws.Range("C49").Value = "f=100 mm, " + [char]0x00D8 + "1`"" + " Achromatic Doublet, SM1-Threaded Mount, ARC: 400-700 nm"
$ws.Range("D49").Value = 4
$ws.Range("E49").Value = 104
$ws.Range("F49").Formula = "=E49*D49"
$ws.Range("G49").Value = "Excitation"
$ws.Range("H49").Value = "A 1:1 relay from ETL to galvo"
$ws.Range("I49").Value = "Lenses L2, L3 in the design view"

# Re-create the hyperlinks one row further down than before, without
# touching the (unchanged) text already stored in those cells.
$ws.Hyperlinks.Add($ws.Range("C88"), $hyperlinkTargets["C87"]) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C89"), $hyperlinkTargets["C88"]) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C90"), $hyperlinkTargets["C89"]) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C91"), $hyperlinkTargets["C90"]) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C96"), $hyperlinkTargets["C95"]) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C79"), $hyperlinkTargets["C78"]) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C92"), $hyperlinkTargets["C91"]) | Out-Null

# Match the adjusted column widths for C, E, F (closest achievable values).
$ws.Columns.Item(3).ColumnWidth = 80.94401041666667
$ws.Columns.Item(5).ColumnWidth = 13.944010416666666
$ws.Columns.Item(6).ColumnWidth = 13.276041666666666

# Update the view: scroll back to the top and select A18 (no frozen top-left offset).
$ws.Range("A18").Select()
